$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values (price cells whose text looks like a plain number
# are first switched to Text format so Excel keeps them as strings, matching
# the original inlineStr cells rather than auto-converting to numbers)
$ws.Range("D2").Value = '37.167.14'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.073.21'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.43'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.676'
$ws.Range("E6").Value = '  +3.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.82'
$ws.Range("E7").Value = '  +21.46%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.97'
$ws.Range("E9").Value = '  +2.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.384'
$ws.Range("E10").Value = '  +4.50%  '
$ws.Range("E11").Value = '  +9.77%  '
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("E13").Value = '  +4.68%  '
$ws.Range("D14").Value = '2.379.07'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.822'
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.48'
$ws.Range("E16").Value = '  +8.29%  '
$ws.Range("D17").Value = '2.082.60'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '37.169.42'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.86'
$ws.Range("E19").Value = '  +3.91%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.33'
$ws.Range("E20").Value = '  +15.70%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0924'
$ws.Range("E21").Value = '  +12.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.51'
$ws.Range("E22").Value = '  +6.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.02'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.99'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.18'
$ws.Range("E27").Value = '  +9.30%  '
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.37'
$ws.Range("E29").Value = '  -1.31%  '
$ws.Range("E30").Value = '  +3.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.76'
$ws.Range("E31").Value = '  +6.21%  '
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0636'
$ws.Range("E33").Value = '  +5.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.42'
$ws.Range("E34").Value = '  +8.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0899'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.30'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("E39").Value = '  +23.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.36'
$ws.Range("E40").Value = '  +2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.24'
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("E43").Value = '  +1.21%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.48'
$ws.Range("E44").Value = '  +26.30%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.55'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.82'
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.68'
$ws.Range("E47").Value = '  +16.19%  '
$ws.Range("E48").Value = '  +10.44%  '
$ws.Range("D49").Value = '1.309.57'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.95'
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.94'
$ws.Range("E51").Value = '  +0.73%  '
